# Update TPM-derived values for Cxadr-Cxadr LR-pair sheet (new TPM recompute).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 3.197839
    "H2" = 9.593517
    "I2" = 0.4734596327704848
    "J2" = 0.4734596327704849
    "M2" = 3.197839
    "N2" = 9.593517
    "O2" = 0.4734596327704848
    "P2" = 0.4734596327704849
    "Q2" = 10.226174269921
    "R2" = 92.035568429289
    "S2" = 0.2241640238631623
    "T2" = 0.2241640238631624
    "G3" = 3.197839
    "H3" = 9.593517
    "I3" = 0.4734596327704848
    "J3" = 0.4734596327704849
    "O3" = 0.2979500593877771
    "P3" = 0.2979500593877771
    "Q3" = 6.435372775507
    "R3" = 57.918354979563
    "S3" = 0.1410673257016811
    "T3" = 0.1410673257016811
    "G4" = 3.197839
    "H4" = 9.593517
    "I4" = 0.4734596327704848
    "J4" = 0.4734596327704849
    "K4" = 1
    "L4" = 0.3333333333333333
    "M4" = 0.1668016666666667
    "N4" = 0.500405
    "O4" = 0.02469600747426772
    "P4" = 0.02469600747426772
    "Q4" = 0.5334048749316667
    "R4" = 4.800643874385
    "S4" = 0.01169256262966394
    "T4" = 0.01169256262966394
    "G5" = 3.197839
    "H5" = 9.593517
    "I5" = 0.4734596327704848
    "J5" = 0.4734596327704849
    "M5" = 1.377142
    "N5" = 4.131426
    "O5" = 0.2038943003674704
    "P5" = 0.2038943003674704
    "Q5" = 4.403878396138
    "R5" = 39.634905565242
    "S5" = 0.09653572057597744
    "T5" = 0.09653572057597745
    "I6" = 0.2979500593877771
    "J6" = 0.2979500593877771
    "M6" = 3.197839
    "N6" = 9.593517
    "O6" = 0.4734596327704848
    "P6" = 0.4734596327704849
    "Q6" = 6.435372775507
    "R6" = 57.918354979563
    "S6" = 0.1410673257016811
    "T6" = 0.1410673257016811
    "I7" = 0.2979500593877771
    "J7" = 0.2979500593877771
    "O7" = 0.2979500593877771
    "P7" = 0.2979500593877771
    "S7" = 0.08877423788917989
    "T7" = 0.08877423788917989
    "I8" = 0.2979500593877771
    "J8" = 0.2979500593877771
    "K8" = 1
    "L8" = 0.3333333333333333
    "M8" = 0.1668016666666667
    "N8" = 0.500405
    "O8" = 0.02469600747426772
    "P8" = 0.02469600747426772
    "Q8" = 0.3356738424216666
    "R8" = 3.021064581795
    "S8" = 0.007358176893599054
    "T8" = 0.007358176893599054
    "I9" = 0.2979500593877771
    "J9" = 0.2979500593877771
    "M9" = 1.377142
    "N9" = 4.131426
    "O9" = 0.2038943003674704
    "P9" = 0.2038943003674704
    "Q9" = 2.771378463646
    "R9" = 24.942406172814
    "S9" = 0.06075031890331705
    "T9" = 0.06075031890331705
    "E10" = 1
    "F10" = 0.3333333333333333
    "G10" = 0.1668016666666667
    "H10" = 0.500405
    "I10" = 0.02469600747426772
    "J10" = 0.02469600747426772
    "M10" = 3.197839
    "N10" = 9.593517
    "O10" = 0.4734596327704848
    "P10" = 0.4734596327704849
    "Q10" = 0.5334048749316667
    "R10" = 4.800643874385
    "S10" = 0.01169256262966394
    "T10" = 0.01169256262966394
    "E11" = 1
    "F11" = 0.3333333333333333
    "G11" = 0.1668016666666667
    "H11" = 0.500405
    "I11" = 0.02469600747426772
    "J11" = 0.02469600747426772
    "O11" = 0.2979500593877771
    "P11" = 0.2979500593877771
    "Q11" = 0.3356738424216666
    "R11" = 3.021064581795
    "S11" = 0.007358176893599054
    "T11" = 0.007358176893599054
    "E12" = 1
    "F12" = 0.3333333333333333
    "G12" = 0.1668016666666667
    "H12" = 0.500405
    "I12" = 0.02469600747426772
    "J12" = 0.02469600747426772
    "K12" = 1
    "L12" = 0.3333333333333333
    "M12" = 0.1668016666666667
    "N12" = 0.500405
    "O12" = 0.02469600747426772
    "P12" = 0.02469600747426772
    "Q12" = 0.02782279600277777
    "R12" = 0.250405164025
    "S12" = 0.0006098927851690872
    "T12" = 0.0006098927851690872
    "E13" = 1
    "F13" = 0.3333333333333333
    "G13" = 0.1668016666666667
    "H13" = 0.500405
    "I13" = 0.02469600747426772
    "J13" = 0.02469600747426772
    "M13" = 1.377142
    "N13" = 4.131426
    "O13" = 0.2038943003674704
    "P13" = 0.2038943003674704
    "Q13" = 0.2297095808366667
    "R13" = 2.06738622753
    "S13" = 0.005035375165835636
    "T13" = 0.005035375165835636
    "G14" = 1.377142
    "H14" = 4.131426
    "I14" = 0.2038943003674704
    "J14" = 0.2038943003674704
    "M14" = 3.197839
    "N14" = 9.593517
    "O14" = 0.4734596327704848
    "P14" = 0.4734596327704849
    "Q14" = 4.403878396138
    "R14" = 39.634905565242
    "S14" = 0.09653572057597744
    "T14" = 0.09653572057597745
    "G15" = 1.377142
    "H15" = 4.131426
    "I15" = 0.2038943003674704
    "J15" = 0.2038943003674704
    "O15" = 0.2979500593877771
    "P15" = 0.2979500593877771
    "Q15" = 2.771378463646
    "R15" = 24.942406172814
    "S15" = 0.06075031890331705
    "T15" = 0.06075031890331705
    "G16" = 1.377142
    "H16" = 4.131426
    "I16" = 0.2038943003674704
    "J16" = 0.2038943003674704
    "K16" = 1
    "L16" = 0.3333333333333333
    "M16" = 0.1668016666666667
    "N16" = 0.500405
    "O16" = 0.02469600747426772
    "P16" = 0.02469600747426772
    "Q16" = 0.2297095808366667
    "R16" = 2.06738622753
    "S16" = 0.005035375165835636
    "T16" = 0.005035375165835636
    "G17" = 1.377142
    "H17" = 4.131426
    "I17" = 0.2038943003674704
    "J17" = 0.2038943003674704
    "M17" = 1.377142
    "N17" = 4.131426
    "O17" = 0.2038943003674704
    "P17" = 0.2038943003674704
    "Q17" = 1.896520088164
    "R17" = 17.068680793476
    "S17" = 0.04157288572234023
    "T17" = 0.04157288572234023
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
